$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed dataset has fewer clusters (57 unique strings vs 69, 56 rows vs 68),
# so drop the now-unused trailing rows first (shifts nothing below them up, since
# they are already the last rows on the sheet).
$ws.Range("A57:B68").EntireRow.Delete()

# Updated cluster names, in row order for A2:A56.
$clusters = @(
    '3323 Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North',
    '3364 Assisi Centre Aged Care Rosanna',
    '3376 Royal Freemasons Coppin Centre Melbourne',
    '3653 Fronditha Thalpori St Albans Aged Care',
    '3825 TLC Forest Lodge Residential Aged Care Frankston North',
    '4167 Royal Freemasons Centennial Lodge Wantirna South',
    '44226 Boneo Primary School Boneo',
    '44321 Maiden Gully Primary School Maiden Gully',
    '45573 Narre Warren South P-12 College Narre Warren South',
    '45695 Sacred Heart Primary School Yarrawonga',
    '46322 Minaret College Officer Campus Officer',
    '50516 Ilim College Glenroy Campus Hadfield',
    '50567 Alamanda K9 College Point Cook',
    '52380 Al Iman College Melton South',
    '52912 Edgars Creek Primary School Wollert',
    '52985 Minaret College Springvale',
    'Adass Israel School Elsternwick',
    'Antonine College Cedar Campus Coburg',
    'Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh',
    'Covenant College Bell Post Hill',
    'Creekside K-9 College Caroline Springs',
    'Dandenong South Primary School Dandenong',
    'Darul Ulum College of Victoria Fawkner October',
    'Derrimut Primary School Derrimut',
    'Devon Meadows Primary School Devon Meadows',
    'Exford Primary School Exford',
    'Flemington Racecourse Flemington',
    'Gilly''s Early Learning Centre Balaclava',
    'Hazel Glen College Doreen',
    'Hazelwood North Primary School Hazelwood North',
    'Ilim College Dallas Main Campus Dallas Oct',
    'Islamic College of Melbourne Tarneit Oct Nov',
    'Lyndhurst Primary School Lyndhurst',
    'Middle Park Primary School Middle Park',
    'Morwell Park Primary School Morwell',
    'Nio Early Learning Adventures Preston',
    'Pentland Primary School Darley',
    'Rutherglen Motor Inn and Walkabout Motel Rutherglen',
    'Sirius College Ibrahim Dellal Campus Sunshine',
    'Sirius College Shepparton Campus Shepparton',
    'Social Gathering Woodvale 30 Oct',
    'Society Restaurant Melbourne',
    'St Ambrose Parish Primary School Woodend',
    'St Joseph''s School Quarry Hill',
    'St Paul''s Primary School Sunshine West',
    'Stevensville Primary School St Albans',
    'Stockdale Road Primary School Traralgon',
    'Supreme Caravans Manufacturing Campbellfield',
    'Templestowe Park Primary School Templestowe',
    'The Lake Primary School Cabarita',
    'Top Yard Rooftop Melbourne',
    'Truganina P-9 College Truganina',
    'Tucker Road Bentleigh Primary School Bentleigh',
    'Warragul Regional College Warragul',
    'Yeshivah College St Kilda East'
)

# Updated active-case counts, in row order for B2:B56.
$cases = @(
    11, 11, 31, 18, 16, 23, 10, 10, 17, 42, 38, 11, 11, 11, 11, 16, 14, 10, 16, 24, 13, 14,
    11, 10, 10, 15, 13, 11, 17, 18, 14, 55, 11, 12, 56, 18, 13, 17, 15, 17, 10, 34, 12, 13,
    18, 10, 18, 51, 20, 12, 14, 12, 11, 20, 24
)

for ($i = 0; $i -lt $clusters.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $clusters[$i]
    $ws.Cells.Item($row, 2).Value = $cases[$i]
}

Write-Host "Done. Used range: $($ws.UsedRange.Address())"
